$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 20-45: existing rows get updated values in columns D,H,I,J,K,L,M,N,P,Q
# (the remaining weekly price observations, columns A,B,C,E,F,G,O,R unchanged).

$ws.Cells.Item(20, 4).Value = 44571
$ws.Cells.Item(20, 8).Value = "Calameño"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 80
$ws.Cells.Item(20, 11).Value = 8000
$ws.Cells.Item(20, 12).Value = 9000
$ws.Cells.Item(20, 13).Value = 8500
$ws.Cells.Item(20, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(20, 16).Value = 472
$ws.Cells.Item(20, 17).Value = 18

$ws.Cells.Item(21, 4).Value = 44571
$ws.Cells.Item(21, 8).Value = "Calameño"
$ws.Cells.Item(21, 9).Value = "Segunda"
$ws.Cells.Item(21, 10).Value = 50
$ws.Cells.Item(21, 11).Value = 5000
$ws.Cells.Item(21, 12).Value = 6000
$ws.Cells.Item(21, 13).Value = 5500
$ws.Cells.Item(21, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(21, 16).Value = 229
$ws.Cells.Item(21, 17).Value = 24

$ws.Cells.Item(22, 4).Value = 44214
$ws.Cells.Item(22, 8).Value = "Calameño"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 160
$ws.Cells.Item(22, 11).Value = 7000
$ws.Cells.Item(22, 12).Value = 8000
$ws.Cells.Item(22, 13).Value = 7500
$ws.Cells.Item(22, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(22, 16).Value = 417
$ws.Cells.Item(22, 17).Value = 18

$ws.Cells.Item(23, 4).Value = 44243
$ws.Cells.Item(23, 8).Value = "Calameño"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 120
$ws.Cells.Item(23, 11).Value = 6500
$ws.Cells.Item(23, 12).Value = 7000
$ws.Cells.Item(23, 13).Value = 6750
$ws.Cells.Item(23, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(23, 16).Value = 375
$ws.Cells.Item(23, 17).Value = 18

$ws.Cells.Item(24, 4).Value = 44243
$ws.Cells.Item(24, 8).Value = "Calameño"
$ws.Cells.Item(24, 9).Value = "Segunda"
$ws.Cells.Item(24, 10).Value = 120
$ws.Cells.Item(24, 11).Value = 5500
$ws.Cells.Item(24, 12).Value = 6000
$ws.Cells.Item(24, 13).Value = 5750
$ws.Cells.Item(24, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(24, 16).Value = 240
$ws.Cells.Item(24, 17).Value = 24

$ws.Cells.Item(25, 4).Value = 44566
$ws.Cells.Item(25, 8).Value = "Calameño"
$ws.Cells.Item(25, 9).Value = "Segunda"
$ws.Cells.Item(25, 10).Value = 50
$ws.Cells.Item(25, 11).Value = 6000
$ws.Cells.Item(25, 12).Value = 7000
$ws.Cells.Item(25, 13).Value = 6500
$ws.Cells.Item(25, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(25, 16).Value = 271
$ws.Cells.Item(25, 17).Value = 24

$ws.Cells.Item(26, 4).Value = 44277
$ws.Cells.Item(26, 8).Value = "Calameño"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 120
$ws.Cells.Item(26, 11).Value = 14000
$ws.Cells.Item(26, 12).Value = 15000
$ws.Cells.Item(26, 13).Value = 14500
$ws.Cells.Item(26, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(26, 16).Value = 806
$ws.Cells.Item(26, 17).Value = 18

$ws.Cells.Item(27, 4).Value = 44277
$ws.Cells.Item(27, 8).Value = "Calameño"
$ws.Cells.Item(27, 9).Value = "Segunda"
$ws.Cells.Item(27, 10).Value = 120
$ws.Cells.Item(27, 11).Value = 13000
$ws.Cells.Item(27, 12).Value = 14000
$ws.Cells.Item(27, 13).Value = 13500
$ws.Cells.Item(27, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(27, 16).Value = 562
$ws.Cells.Item(27, 17).Value = 24

$ws.Cells.Item(28, 4).Value = 44181
$ws.Cells.Item(28, 8).Value = "Tuna"
$ws.Cells.Item(28, 9).Value = "Segunda"
$ws.Cells.Item(28, 10).Value = 100
$ws.Cells.Item(28, 11).Value = 10000
$ws.Cells.Item(28, 12).Value = 12000
$ws.Cells.Item(28, 13).Value = 11000
$ws.Cells.Item(28, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(28, 16).Value = 458
$ws.Cells.Item(28, 17).Value = 24

$ws.Cells.Item(29, 4).Value = 44567
$ws.Cells.Item(29, 8).Value = "Tuna"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 60
$ws.Cells.Item(29, 11).Value = 8000
$ws.Cells.Item(29, 12).Value = 9000
$ws.Cells.Item(29, 13).Value = 8500
$ws.Cells.Item(29, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(29, 16).Value = 472
$ws.Cells.Item(29, 17).Value = 18

$ws.Cells.Item(30, 4).Value = 44176
$ws.Cells.Item(30, 8).Value = "Tuna"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 140
$ws.Cells.Item(30, 11).Value = 10000
$ws.Cells.Item(30, 12).Value = 11000
$ws.Cells.Item(30, 13).Value = 10500
$ws.Cells.Item(30, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(30, 16).Value = 583
$ws.Cells.Item(30, 17).Value = 18

$ws.Cells.Item(31, 4).Value = 44218
$ws.Cells.Item(31, 8).Value = "Calameño"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 80
$ws.Cells.Item(31, 11).Value = 9000
$ws.Cells.Item(31, 12).Value = 10000
$ws.Cells.Item(31, 13).Value = 9500
$ws.Cells.Item(31, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(31, 16).Value = 528
$ws.Cells.Item(31, 17).Value = 18

$ws.Cells.Item(32, 4).Value = 44218
$ws.Cells.Item(32, 8).Value = "Calameño"
$ws.Cells.Item(32, 9).Value = "Segunda"
$ws.Cells.Item(32, 10).Value = 120
$ws.Cells.Item(32, 11).Value = 8000
$ws.Cells.Item(32, 12).Value = 9000
$ws.Cells.Item(32, 13).Value = 8500
$ws.Cells.Item(32, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(32, 16).Value = 354
$ws.Cells.Item(32, 17).Value = 24

$ws.Cells.Item(33, 4).Value = 44218
$ws.Cells.Item(33, 8).Value = "Tuna"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 80
$ws.Cells.Item(33, 11).Value = 11000
$ws.Cells.Item(33, 12).Value = 12000
$ws.Cells.Item(33, 13).Value = 11500
$ws.Cells.Item(33, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(33, 16).Value = 639
$ws.Cells.Item(33, 17).Value = 18

$ws.Cells.Item(34, 4).Value = 44218
$ws.Cells.Item(34, 8).Value = "Tuna"
$ws.Cells.Item(34, 9).Value = "Segunda"
$ws.Cells.Item(34, 10).Value = 120
$ws.Cells.Item(34, 11).Value = 8000
$ws.Cells.Item(34, 12).Value = 9000
$ws.Cells.Item(34, 13).Value = 8500
$ws.Cells.Item(34, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(34, 16).Value = 354
$ws.Cells.Item(34, 17).Value = 24

$ws.Cells.Item(35, 4).Value = 44259
$ws.Cells.Item(35, 8).Value = "Calameño"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 80
$ws.Cells.Item(35, 11).Value = 7500
$ws.Cells.Item(35, 12).Value = 8000
$ws.Cells.Item(35, 13).Value = 7750
$ws.Cells.Item(35, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(35, 16).Value = 431
$ws.Cells.Item(35, 17).Value = 18

$ws.Cells.Item(36, 4).Value = 44259
$ws.Cells.Item(36, 8).Value = "Calameño"
$ws.Cells.Item(36, 9).Value = "Segunda"
$ws.Cells.Item(36, 10).Value = 100
$ws.Cells.Item(36, 11).Value = 6500
$ws.Cells.Item(36, 12).Value = 7000
$ws.Cells.Item(36, 13).Value = 6750
$ws.Cells.Item(36, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(36, 16).Value = 281
$ws.Cells.Item(36, 17).Value = 24

$ws.Cells.Item(37, 4).Value = 44251
$ws.Cells.Item(37, 8).Value = "Tuna"
$ws.Cells.Item(37, 9).Value = "Segunda"
$ws.Cells.Item(37, 10).Value = 120
$ws.Cells.Item(37, 11).Value = 7000
$ws.Cells.Item(37, 12).Value = 8000
$ws.Cells.Item(37, 13).Value = 7500
$ws.Cells.Item(37, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(37, 16).Value = 312
$ws.Cells.Item(37, 17).Value = 24

$ws.Cells.Item(38, 4).Value = 44559
$ws.Cells.Item(38, 8).Value = "Calameño"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 80
$ws.Cells.Item(38, 11).Value = 8000
$ws.Cells.Item(38, 12).Value = 9000
$ws.Cells.Item(38, 13).Value = 8500
$ws.Cells.Item(38, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(38, 16).Value = 531
$ws.Cells.Item(38, 17).Value = 16

$ws.Cells.Item(39, 4).Value = 44559
$ws.Cells.Item(39, 8).Value = "Calameño"
$ws.Cells.Item(39, 9).Value = "Super"
$ws.Cells.Item(39, 10).Value = 70
$ws.Cells.Item(39, 11).Value = 12000
$ws.Cells.Item(39, 12).Value = 13000
$ws.Cells.Item(39, 13).Value = 12500
$ws.Cells.Item(39, 14).Value = "$/caja 12 unidades"
$ws.Cells.Item(39, 16).Value = 1042
$ws.Cells.Item(39, 17).Value = 12

$ws.Cells.Item(40, 4).Value = 44208
$ws.Cells.Item(40, 8).Value = "Calameño"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 140
$ws.Cells.Item(40, 11).Value = 9000
$ws.Cells.Item(40, 12).Value = 10000
$ws.Cells.Item(40, 13).Value = 9500
$ws.Cells.Item(40, 14).Value = "$/caja 12 unidades"
$ws.Cells.Item(40, 16).Value = 792
$ws.Cells.Item(40, 17).Value = 12

$ws.Cells.Item(41, 4).Value = 44208
$ws.Cells.Item(41, 8).Value = "Tuna"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 120
$ws.Cells.Item(41, 11).Value = 8000
$ws.Cells.Item(41, 12).Value = 9000
$ws.Cells.Item(41, 13).Value = 8500
$ws.Cells.Item(41, 14).Value = "$/caja 14 unidades"
$ws.Cells.Item(41, 16).Value = 607
$ws.Cells.Item(41, 17).Value = 14

$ws.Cells.Item(42, 4).Value = 44200
$ws.Cells.Item(42, 8).Value = "Calameño"
$ws.Cells.Item(42, 9).Value = "Extra (muy buena)"
$ws.Cells.Item(42, 10).Value = 100
$ws.Cells.Item(42, 11).Value = 9000
$ws.Cells.Item(42, 12).Value = 9500
$ws.Cells.Item(42, 13).Value = 9250
$ws.Cells.Item(42, 14).Value = "$/caja 12 unidades"
$ws.Cells.Item(42, 16).Value = 771
$ws.Cells.Item(42, 17).Value = 12

$ws.Cells.Item(43, 4).Value = 44200
$ws.Cells.Item(43, 8).Value = "Tuna"
$ws.Cells.Item(43, 9).Value = "Segunda"
$ws.Cells.Item(43, 10).Value = 120
$ws.Cells.Item(43, 11).Value = 6500
$ws.Cells.Item(43, 12).Value = 7000
$ws.Cells.Item(43, 13).Value = 6750
$ws.Cells.Item(43, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(43, 16).Value = 281
$ws.Cells.Item(43, 17).Value = 24

$ws.Cells.Item(44, 4).Value = 44224
$ws.Cells.Item(44, 8).Value = "Calameño"
$ws.Cells.Item(44, 9).Value = "Segunda"
$ws.Cells.Item(44, 10).Value = 120
$ws.Cells.Item(44, 11).Value = 5000
$ws.Cells.Item(44, 12).Value = 6000
$ws.Cells.Item(44, 13).Value = 5500
$ws.Cells.Item(44, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(44, 16).Value = 229
$ws.Cells.Item(44, 17).Value = 24

$ws.Cells.Item(45, 4).Value = 44568
$ws.Cells.Item(45, 8).Value = "Calameño"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 70
$ws.Cells.Item(45, 11).Value = 9000
$ws.Cells.Item(45, 12).Value = 10000
$ws.Cells.Item(45, 13).Value = 9500
$ws.Cells.Item(45, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(45, 16).Value = 528
$ws.Cells.Item(45, 17).Value = 18

# New rows 46 and 47: append two more weekly price observations.

$ws.Cells.Item(46, 1).Value = 1
$ws.Cells.Item(46, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(46, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(46, 4).Value = 44568
$ws.Cells.Item(46, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(46, 5).Value = 15
$ws.Cells.Item(46, 6).Value = 100112027
$ws.Cells.Item(46, 7).Value = "Melón"
$ws.Cells.Item(46, 8).Value = "Calameño"
$ws.Cells.Item(46, 9).Value = "Segunda"
$ws.Cells.Item(46, 10).Value = 100
$ws.Cells.Item(46, 11).Value = 5000
$ws.Cells.Item(46, 12).Value = 5500
$ws.Cells.Item(46, 13).Value = 5250
$ws.Cells.Item(46, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(46, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(46, 16).Value = 219
$ws.Cells.Item(46, 17).Value = 24
$ws.Cells.Item(46, 18).Value = "Hortaliza"

$ws.Cells.Item(47, 1).Value = 1
$ws.Cells.Item(47, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(47, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(47, 4).Value = 44568
$ws.Cells.Item(47, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(47, 5).Value = 15
$ws.Cells.Item(47, 6).Value = 100112027
$ws.Cells.Item(47, 7).Value = "Melón"
$ws.Cells.Item(47, 8).Value = "Tuna"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 80
$ws.Cells.Item(47, 11).Value = 10000
$ws.Cells.Item(47, 12).Value = 11000
$ws.Cells.Item(47, 13).Value = 10500
$ws.Cells.Item(47, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(47, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(47, 16).Value = 583
$ws.Cells.Item(47, 17).Value = 18
$ws.Cells.Item(47, 18).Value = "Hortaliza"
